$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Insert a new row at position 73 (shifts rows 73.. down by one)
$ws.Rows("73:73").Insert()

# Match the formatting of the surrounding rows (Insert() alone does not
# copy the row-above formatting the way interactive Excel does)
$ws.Range("A74:H74").Copy()
$ws.Range("A73:H73").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new row 73 with the "neovim" entry (installed via git-bash / Work column)
$ws.Range("A73").Value = "neovim"
$ws.Range("D73").Value = "x"

# Mark "x" for fzf (row 33) in the Work (git-bash) column
$ws.Range("D33").Value = "x"

# Mark "x" for ripgrep (now row 89 after the insert) in the Work (git-bash) column
$ws.Range("D89").Value = "x"

# Re-apply the AutoFilter over the new, larger data range (A1:H130)
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}
$ws.Range("A1:H130").AutoFilter()

# Keep the hidden _FilterDatabase defined name in sync with the new range
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Tabelle1!_FilterDatabase") {
        $n.RefersTo = "=Tabelle1!`$A`$1:`$H`$130"
    }
}

# Update the view: scroll position and active cell selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 65
$ws.Range("E73").Select()
